# edit.ps1
# Applies the "feat: add 2022-Q3 data" change:
#  - Inserts a new worksheet "2022-Q3" right after "总计", pushing all the
#    other quarterly sheets one position to the right (their own data is
#    untouched).
#  - Populates "2022-Q3" with its two fund rows.
#  - Updates the "总计" summary sheet: a new top data row for 2022-Q3 is
#    added and the previously-existing quarters shift down one row, with
#    a brand new last row for 2021-Q1.

$wb = $excel.ActiveWorkbook

function Set-TextCell($sheet, $row, $col, $val) {
    # Forces the cell to be written as text (keeps leading zeros / avoids
    # Excel's automatic number coercion), then strips the temporary
    # number-format override so the cell ends up with the default style,
    # matching cells elsewhere in the workbook that hold text-looking values.
    $c = $sheet.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q3" sheet by duplicating the existing
#    "2022-Q2" sheet (this keeps identical headers/column styles/sheet
#    properties), placing the copy right after "总计".
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)
$wsQ2Source = $wb.Worksheets.Item(2)
$wsQ2Source.Copy([System.Reflection.Missing]::Value, $wsTotal)

$wsQ3 = $wb.Worksheets.Item(2)
$wsQ3.Name = "2022-Q3"

# The duplicated sheet has 5 data rows (rows 2-6); 2022-Q3 only needs 2
# data rows, so remove the extra ones.
$wsQ3.Rows("4:6").Delete()

# Row 2: 002137 / 诺安利鑫灵活配置混合A
$wsQ3.Cells.Item(2,1).Value = 0
Set-TextCell $wsQ3 2 2 "002137"
Set-TextCell $wsQ3 2 3 "诺安利鑫灵活配置混合A"
Set-TextCell $wsQ3 2 4 "0.44"
Set-TextCell $wsQ3 2 5 "76.46"
Set-TextCell $wsQ3 2 6 "3.16"
Set-TextCell $wsQ3 2 7 "0.0139"
$wsQ3.Cells.Item(2,8).Value = 7

# Row 3: 014521 / 诺安利鑫灵活配置混合C
$wsQ3.Cells.Item(3,1).Value = 1
Set-TextCell $wsQ3 3 2 "014521"
Set-TextCell $wsQ3 3 3 "诺安利鑫灵活配置混合C"
Set-TextCell $wsQ3 3 4 "0.05"
Set-TextCell $wsQ3 3 5 "76.46"
Set-TextCell $wsQ3 3 6 "3.16"
Set-TextCell $wsQ3 3 7 "0.0016"
$wsQ3.Cells.Item(3,8).Value = 7

# ---------------------------------------------------------------------
# 2) Update the "总计" (total) sheet: insert 2022-Q3 numbers at the top
#    and shift the rest down by one, appending the new last row for
#    2021-Q1.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)

# Extend formatting of column A down to the brand-new row 7 (matches the
# style already used by A2:A6).
$wsTotal.Range("A6").Copy()
$wsTotal.Range("A7").PasteSpecial(-4122)

$wsTotal.Cells.Item(7,1).Value = 5
$wsTotal.Cells.Item(7,2).Value = "2021-Q1"
$wsTotal.Cells.Item(7,3).Value = 1
$wsTotal.Cells.Item(7,4).Value = 0.5

$wsTotal.Cells.Item(2,2).Value = "2022-Q3"
$wsTotal.Cells.Item(2,3).Value = 2
$wsTotal.Cells.Item(2,4).Value = 0.02

$wsTotal.Cells.Item(3,2).Value = "2022-Q2"
$wsTotal.Cells.Item(3,3).Value = 5
$wsTotal.Cells.Item(3,4).Value = 0.27

$wsTotal.Cells.Item(4,2).Value = "2022-Q1"
$wsTotal.Cells.Item(4,3).Value = 3
$wsTotal.Cells.Item(4,4).Value = 0.24

$wsTotal.Cells.Item(5,2).Value = "2021-Q4"
$wsTotal.Cells.Item(5,3).Value = 3
$wsTotal.Cells.Item(5,4).Value = 0.01

$wsTotal.Cells.Item(6,2).Value = "2021-Q3"
$wsTotal.Cells.Item(6,3).Value = 1
$wsTotal.Cells.Item(6,4).Value = 0
